$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 242 (shifts existing rows 242:307 down to 243:308)
$ws.Rows.Item(242).Insert()

# Populate the newly inserted row with the new weekly price entry
$ws.Range("A242").Value = 3
$ws.Range("B242").Value = "Femacal de La Calera"
$ws.Range("C242").Value = "Coquimbo"
$ws.Range("D242").Value = 44642
$ws.Range("E242").Value = 5
$ws.Range("F242").Value = 100112012
$ws.Range("G242").Value = "Espinaca"
$ws.Range("H242").Value = "Sin especificar"
$ws.Range("I242").Value = "Primera"
$ws.Range("J242").Value = 130
$ws.Range("K242").Value = 4000
$ws.Range("L242").Value = 4500
$ws.Range("M242").Value = 4269
$ws.Range("N242").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O242").Value = "Provincia de Quillota"
$ws.Range("P242").Value = 1423
$ws.Range("Q242").Value = 3
$ws.Range("R242").Value = "Hortaliza"
